# Updated cryptos list on Tue Mar 26 14:09:38 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '70.405.16'
Set-TextValue 'E2' '  +3.57%  '
Set-TextValue 'D3' '3.611.08'
Set-TextValue 'E3' '  +3.55%  '
Set-TextValue 'E4' '  -0.11%  '
Set-TextValue 'D5' '583.10'
Set-TextValue 'E5' '  +0.10%  '
Set-TextValue 'D6' '192.53'
Set-TextValue 'E6' '  +1.32%  '
Set-TextValue 'D7' '0.637'
Set-TextValue 'E7' '  +0.80%  '
Set-TextValue 'D8' '3.609.51'
Set-TextValue 'E8' '  +3.87%  '
Set-TextValue 'D9' '0.999'
Set-TextValue 'E9' '  -0.13%  '
Set-TextValue 'D10' '0.181'
Set-TextValue 'E10' '  +3.76%  '
Set-TextValue 'D11' '0.666'
Set-TextValue 'E11' '  +3.07%  '
Set-TextValue 'D12' '57.13'
Set-TextValue 'E12' '  -1.61%  '
Set-TextValue 'D13' '0.0000304'
Set-TextValue 'E13' '  +8.95%  '
Set-TextValue 'D14' '9.81'
Set-TextValue 'E14' '  +3.47%  '
Set-TextValue 'D15' '4.186.27'
Set-TextValue 'E15' '  +4.08%  '
Set-TextValue 'D16' '20.31'
Set-TextValue 'E16' '  +6.68%  '
Set-TextValue 'D17' '3.612.80'
Set-TextValue 'E17' '  +3.73%  '
Set-TextValue 'D18' '70.334.75'
Set-TextValue 'E18' '  +3.59%  '
Set-TextValue 'D19' '12.62'
Set-TextValue 'E19' '  +3.61%  '
Set-TextValue 'E20' '  +2.51%  '
Set-TextValue 'D21' '1.05'
Set-TextValue 'E21' '  +2.71%  '
Set-TextValue 'D22' '484.59'
Set-TextValue 'E22' '  -0.65%  '
Set-TextValue 'D23' '19.36'
Set-TextValue 'E23' '  +11.60%  '
Set-TextValue 'D24' '5.10'
Set-TextValue 'E24' '  -9.11%  '
Set-TextValue 'D25' '4.40'
Set-TextValue 'E25' '  +0.98%  '
Set-TextValue 'D26' '90.03'
Set-TextValue 'E26' '  +0.27%  '
Set-TextValue 'D27' '3.11'
Set-TextValue 'E27' '  +3.63%  '
Set-TextValue 'D28' '11.27'
Set-TextValue 'E28' '  +2.67%  '
Set-TextValue 'D29' '9.42'
Set-TextValue 'E29' '  +3.49%  '
Set-TextValue 'B30' 'EthereumClassic'
Set-TextValue 'C30' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D30' '32.52'
Set-TextValue 'E30' '  +3.47%  '
Set-TextValue 'B31' 'NEARProtocol'
Set-TextValue 'C31' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D31' '7.80'
Set-TextValue 'E31' '  +5.61%  '
Set-TextValue 'E32' '  +10.60%  '
Set-TextValue 'D33' '12.20'
Set-TextValue 'E33' '  +2.82%  '
Set-TextValue 'D34' '66.48'
Set-TextValue 'E34' '  +2.89%  '
Set-TextValue 'D35' '605.84'
Set-TextValue 'E35' '  -0.39%  '
Set-TextValue 'D36' '39.89'
Set-TextValue 'E36' '  +7.45%  '
Set-TextValue 'D37' '0.0₃0827'
Set-TextValue 'E37' '  +6.61%  '
Set-TextValue 'D38' '0.405'
Set-TextValue 'E38' '  +4.12%  '
Set-TextValue 'D39' '0.147'
Set-TextValue 'E39' '  +0.15%  '
Set-TextValue 'D40' '0.998'
Set-TextValue 'E40' '  -0.24%  '
Set-TextValue 'D41' '2.97'
Set-TextValue 'E41' '  +14.75%  '
Set-TextValue 'D42' '3.54'
Set-TextValue 'E42' '  +1.30%  '
Set-TextValue 'D43' '3.309.33'
Set-TextValue 'E43' '  +2.31%  '
Set-TextValue 'D44' '3.20'
Set-TextValue 'E44' '  +18.82%  '
Set-TextValue 'D45' '3.12'
Set-TextValue 'E45' '  +7.32%  '
Set-TextValue 'D46' '0.0452'
Set-TextValue 'E46' '  +4.26%  '
Set-TextValue 'D47' '9.62'
Set-TextValue 'E47' '  +10.42%  '
Set-TextValue 'D48' '3.37'
Set-TextValue 'E48' '  +4.12%  '
Set-TextValue 'D49' '0.138'
Set-TextValue 'E49' '  +1.91%  '
Set-TextValue 'D50' '1.00'
Set-TextValue 'E50' '  +0.11%  '
Set-TextValue 'E51' '  +1.50%  '
